$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record at row 5 (pushes the existing rows 5-29 down to 6-30,
# preserving their values/formatting automatically).
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5,3).Value = "Los Lagos"
$ws.Cells.Item(5,4).Value = 44831
$ws.Cells.Item(5,5).Value = 10
$ws.Cells.Item(5,6).Value = 100112013
$ws.Cells.Item(5,7).Value = "Alcachofa"
$ws.Cells.Item(5,8).Value = "Madrigal"
$ws.Cells.Item(5,9).Value = "Primera"
$ws.Cells.Item(5,10).Value = 180
$ws.Cells.Item(5,11).Value = 12000
$ws.Cells.Item(5,12).Value = 13000
$ws.Cells.Item(5,13).Value = 12444
$ws.Cells.Item(5,14).Value = "$/caja 40 unidades"
$ws.Cells.Item(5,15).Value = "Provincia de Limarí"
$ws.Cells.Item(5,16).Value = 311
$ws.Cells.Item(5,17).Value = 40
$ws.Cells.Item(5,18).Value = "Hortaliza"

# Append a new weekly record as the new last row (31).
$ws.Cells.Item(31,1).Value = 4
$ws.Cells.Item(31,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(31,3).Value = "Los Lagos"
$ws.Cells.Item(31,4).Value = 44832
$ws.Cells.Item(31,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(31,5).Value = 10
$ws.Cells.Item(31,6).Value = 100112013
$ws.Cells.Item(31,7).Value = "Alcachofa"
$ws.Cells.Item(31,8).Value = "Madrigal"
$ws.Cells.Item(31,9).Value = "Primera"
$ws.Cells.Item(31,10).Value = 120
$ws.Cells.Item(31,11).Value = 12000
$ws.Cells.Item(31,12).Value = 13000
$ws.Cells.Item(31,13).Value = 12500
$ws.Cells.Item(31,14).Value = "$/caja 40 unidades"
$ws.Cells.Item(31,15).Value = "Provincia de Limarí"
$ws.Cells.Item(31,16).Value = 312
$ws.Cells.Item(31,17).Value = 40
$ws.Cells.Item(31,18).Value = "Hortaliza"
